$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.167.65'
$ws.Range('E2').Value = '  +0.55%  '

$ws.Range('D3').Value = '2.637.53'
$ws.Range('E3').Value = '  +0.48%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').Value = '598.47'
$ws.Range('E5').Value = '  +0.59%  '

$ws.Range('D6').Value = '154.46'
$ws.Range('E6').Value = '  +0.71%  '

$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('D8').Value = '0.543'
$ws.Range('E8').Value = '  -0.68%  '

$ws.Range('D9').Value = '2.636.68'
$ws.Range('E9').Value = '  +0.51%  '

$ws.Range('D10').Value = '0.143'
$ws.Range('E10').Value = '  +7.29%  '

$ws.Range('E11').Value = '  -0.72%  '

$ws.Range('E12').Value = '  +0.29%  '

$ws.Range('D13').Value = '0.352'
$ws.Range('E13').Value = '  +1.54%  '

$ws.Range('D14').Value = '27.80'
$ws.Range('E14').Value = '  +0.90%  '

$ws.Range('D15').Value = '0.0000191'
$ws.Range('E15').Value = '  +2.28%  '

$ws.Range('D16').Value = '3.118.44'
$ws.Range('E16').Value = '  +0.53%  '

$ws.Range('D17').Value = '68.065.94'
$ws.Range('E17').Value = '  +0.53%  '

$ws.Range('D18').Value = '2.636.84'
$ws.Range('E18').Value = '  -0.13%  '

$ws.Range('D19').Value = '11.33'
$ws.Range('E19').Value = '  -0.37%  '

$ws.Range('D20').Value = '362.50'
$ws.Range('E20').Value = '  -1.15%  '

$ws.Range('D21').Value = '7.40'
$ws.Range('E21').Value = '  -0.23%  '

$ws.Range('D22').Value = '4.34'

$ws.Range('E23').Value = '  +0.72%  '

$ws.Range('D24').Value = '2.05'
$ws.Range('E24').Value = '  -1.48%  '

$ws.Range('D25').Value = '75.22'
$ws.Range('E25').Value = '  +4.54%  '

$ws.Range('E26').Value = '  +0.04%  '

$ws.Range('D27').Value = '9.74'
$ws.Range('E27').Value = '  -0.79%  '

$ws.Range('E28').Value = '  +1.61%  '

$ws.Range('D29').Value = '2.774.93'
$ws.Range('E29').Value = '  +0.53%  '

$ws.Range('E30').Value = '  +0.00%  '

$ws.Range('D31').Value = '559.58'
$ws.Range('E31').Value = '  -2.81%  '

$ws.Range('D32').Value = '8.01'
$ws.Range('E32').Value = '  +1.28%  '

$ws.Range('D33').Value = '1.39'
$ws.Range('E33').Value = '  -0.03%  '

$ws.Range('E34').Value = '  +0.94%  '

$ws.Range('E35').Value = '  +1.24%  '

$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  -0.04%  '

$ws.Range('D37').Value = '1.57'
$ws.Range('E37').Value = '  +2.44%  '

$ws.Range('D38').Value = '160.79'
$ws.Range('E38').Value = '  +0.72%  '

$ws.Range('D39').Value = '19.26'
$ws.Range('E39').Value = '  +0.81%  '

$ws.Range('D40').Value = '0.372'
$ws.Range('E40').Value = '  +1.53%  '

$ws.Range('E41').Value = '  -0.11%  '

$ws.Range('D42').Value = '5.31'
$ws.Range('E42').Value = '  -0.59%  '

$ws.Range('D43').Value = '0.0₆0337'
$ws.Range('E43').Value = '  +1.37%  '

$ws.Range('E44').Value = '  +2.61%  '

$ws.Range('E45').Value = '  -0.94%  '

$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  -0.01%  '

$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').Value = '40.60'
$ws.Range('E47').Value = '  +1.24%  '

$ws.Range('D48').Value = '156.95'
$ws.Range('E48').Value = '  +1.27%  '

$ws.Range('D49').Value = '3.73'
$ws.Range('E49').Value = '  +1.47%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.0785'
$ws.Range('E50').Value = '  +0.94%  '

$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '21.73'
$ws.Range('E51').Value = '  -0.37%  '
